# Remove lab modules that don't exist
# Clears the "topic" values for lab modules lab03 (D24), lab04 (D25), and lab07 (D28)
# which referenced labs that don't actually exist.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D24").ClearContents()
$ws.Range("D25").ClearContents()
$ws.Range("D28").ClearContents()

# Update the selection to match the resulting workbook state
$ws.Range("D24:D28").Select()
